# Manually change users and multiple excel writes
# Update the three per-person "Drive Summary" sheets with this period's
# odometer / points figures.

$wb = $excel.ActiveWorkbook

# --- Stefan ---------------------------------------------------------
$ws = $wb.Worksheets.Item("Stefan")
$ws.Range("F6").Value = 828
$ws.Range("G6").Value = 542
$ws.Range("C7").Value = 256
$ws.Range("C9").Select

# --- Christiaan -------------------------------------------------------
$ws = $wb.Worksheets.Item("Christiaan")
$ws.Range("F6").Value = 407
$ws.Range("G6").Value = 547
$ws.Range("C7").Value = 188
$ws.Range("C13").Select

# --- Derrick ----------------------------------------------------------
$ws = $wb.Worksheets.Item("Derrick")
$ws.Range("C6").Value = 150
$ws.Range("F6").Value = 316
$ws.Range("G6").Value = 548
$ws.Range("C7").Value = 300
$ws.Range("C9").Value = 300
$ws.Range("C13").Value = 1179.26
$ws.Range("G6").Select
